# Update "想去人数" (want-to-go count) values in column F across all sheets
# to match the newly generated data snapshot, per commit:
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# Map of sheet name -> (row -> new F value)
$updates = @{
    "展览" = @{
        2  = 321
        3  = 1268
        4  = 364
        5  = 334
        6  = 3845
        8  = 753
        9  = 2231
        11 = 220
        13 = 162
        14 = 163
        15 = 2145
        18 = 48
        19 = 337
        21 = 27
        22 = 269
    }
    "演出" = @{
        3  = 42
        8  = 127
        11 = 84
    }
    "本地生活" = @{
        3 = 821
        4 = 2083
        5 = 315
    }
    "全部类型" = @{
        3  = 821
        4  = 2083
        5  = 315
        7  = 42
        8  = 42
        10 = 321
        11 = 1268
        12 = 364
        16 = 334
        17 = 3845
        18 = 127
        22 = 84
        23 = 753
        24 = 2231
        27 = 220
        29 = 162
        30 = 163
        32 = 2146
        37 = 48
        38 = 337
        40 = 27
        48 = 269
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowMap = $updates[$sheetName]
    foreach ($row in $rowMap.Keys) {
        $ws.Range("F$row").Value = $rowMap[$row]
    }
}
